$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 205759.8

$ws.Range("H53").Value = 296.06897
$ws.Range("I53").Value = 237.8125
$ws.Range("K53").Value = 237.8125
$ws.Range("M53").Value = 399.1875

$ws.Range("H64").Value = 69688.60000000001
$ws.Range("J64").Value = 3369.889
$ws.Range("L64").Value = 3369.889
$ws.Range("N64").Value = -3865.889

$ws.Range("H67").Value = 69688.60000000001
$ws.Range("J67").Value = 3369.889
$ws.Range("L67").Value = 3369.889
$ws.Range("N67").Value = -5085.889

$ws.Range("H70").Value = 1528.7084
$ws.Range("J70").Value = 1555.7142
$ws.Range("L70").Value = 4667.142599999999
$ws.Range("N70").Value = -5207.142599999999

$ws.Range("H73").Value = 1528.7084
$ws.Range("J73").Value = 1555.7142
$ws.Range("L73").Value = 4667.142599999999
$ws.Range("N73").Value = -6539.142599999999

$ws.Range("H76").Value = 4533.154
$ws.Range("I76").Value = 4442.9
$ws.Range("J76").Value = 4834
$ws.Range("K76").Value = 4442.9
$ws.Range("L76").Value = 4834
$ws.Range("M76").Value = -4127.9
$ws.Range("N76").Value = -5464

$ws.Range("H79").Value = 4533.154
$ws.Range("I79").Value = 4442.9
$ws.Range("J79").Value = 4834
$ws.Range("K79").Value = 4442.9
$ws.Range("L79").Value = 4834
$ws.Range("M79").Value = -3350.9
$ws.Range("N79").Value = -7018

$ws.Range("H82").Value = 2451.8462
$ws.Range("I82").Value = 1124.8572
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 3374.5716
$ws.Range("L82").Value = 12000
$ws.Range("M82").Value = -2968.5716
$ws.Range("N82").Value = -12812

$ws.Range("H85").Value = 2451.8462
$ws.Range("I85").Value = 1124.8572
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 3374.5716
$ws.Range("L85").Value = 12000
$ws.Range("M85").Value = -1970.5716
$ws.Range("N85").Value = -14808

$ws.Range("H116").Value = 4500
$ws.Range("I116").Value = 4500
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 4500
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -1058
$ws.Range("N116").ClearContents()

$ws.Range("H118").Value = 6843.75
$ws.Range("I118").Value = 8654.166999999999
$ws.Range("J118").Value = 1412.5
$ws.Range("K118").Value = 25962.501
$ws.Range("L118").Value = 4237.5
$ws.Range("M118").Value = -24305.501
$ws.Range("N118").Value = -7551.5

$ws.Range("H121").Value = 749.28
$ws.Range("I121").Value = 1333.3334
$ws.Range("J121").Value = 669.63635
$ws.Range("K121").Value = 4000.0002
$ws.Range("L121").Value = 2008.90905
$ws.Range("M121").Value = -2253.0002
$ws.Range("N121").Value = -5502.90905

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H138").Value = 2508.62
$ws.Range("I138").Value = 1180.95
$ws.Range("J138").Value = 2840.5376
$ws.Range("K138").Value = 3542.85
$ws.Range("L138").Value = 8521.612800000001
$ws.Range("M138").Value = 1597.15
$ws.Range("N138").Value = -18801.6128

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1812.5834
$ws.Range("I61").Value = 1332.5
$ws.Range("J61").Value = 2292.6667
$ws.Range("K61").Value = 1332.5
$ws.Range("L61").Value = 2292.6667
$ws.Range("M61").Value = -1120.5
$ws.Range("N61").Value = -2716.6667

$ws.Range("H63").Value = 2398
$ws.Range("J63").Value = 2796.6667
$ws.Range("L63").Value = 2796.6667
$ws.Range("N63").Value = -4168.6667

$ws.Range("H66").Value = 2398
$ws.Range("J66").Value = 2796.6667
$ws.Range("L66").Value = 13983.3335
$ws.Range("N66").Value = -20847.3335

$ws.Range("H132").Value = 24299.143
$ws.Range("I132").Value = 37949.35
$ws.Range("K132").Value = 113848.05
$ws.Range("M132").Value = -111318.05

$ws.Range("H136").Value = 1812.5834
$ws.Range("I136").Value = 1332.5
$ws.Range("J136").Value = 2292.6667
$ws.Range("K136").Value = 3997.5
$ws.Range("L136").Value = 6878.000100000001
$ws.Range("M136").Value = -1447.5
$ws.Range("N136").Value = -11978.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17789.889
$ws.Range("I82").Value = 2014
$ws.Range("J82").Value = 30410.6
$ws.Range("K82").Value = 2014
$ws.Range("L82").Value = 30410.6
$ws.Range("M82").Value = -1631
$ws.Range("N82").Value = -31176.6

$ws.Range("H85").Value = 17789.889
$ws.Range("I85").Value = 2014
$ws.Range("J85").Value = 30410.6
$ws.Range("K85").Value = 2014
$ws.Range("L85").Value = 30410.6
$ws.Range("M85").Value = -688
$ws.Range("N85").Value = -33062.6

$ws.Range("H134").Value = 3321.7632
$ws.Range("I134").Value = 3130.7273
$ws.Range("K134").Value = 9392.1819
$ws.Range("M134").Value = -6857.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 3820.8125
$ws.Range("I107").Value = 4646.12
$ws.Range("J107").Value = 873.2857
$ws.Range("K107").Value = 4646.12
$ws.Range("L107").Value = 873.2857
$ws.Range("M107").Value = -2726.12
$ws.Range("N107").Value = -4713.2857

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1142.8572
$ws.Range("I17").Value = 775
$ws.Range("J17").Value = 1633.3334
$ws.Range("K17").Value = 2325
$ws.Range("L17").Value = 4900.0002
$ws.Range("M17").Value = -2156
$ws.Range("N17").Value = -5238.0002

$ws.Range("H131").Value = 539646.75
$ws.Range("I131").Value = 710.0526
$ws.Range("J131").Value = 676177.4
$ws.Range("K131").Value = 2130.1578
$ws.Range("L131").Value = 2028532.2
$ws.Range("M131").Value = 2909.8422
$ws.Range("N131").Value = -2038612.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 163970.92
$ws.Range("I102").Value = 1624.0416
$ws.Range("K102").Value = 1624.0416
$ws.Range("M102").Value = -2.041600000000017

$ws.Range("H122").Value = 2454.8235
$ws.Range("I122").Value = 2066.0908
$ws.Range("J122").Value = 3167.5
$ws.Range("K122").Value = 6198.2724
$ws.Range("L122").Value = 9502.5
$ws.Range("M122").Value = -3748.2724
$ws.Range("N122").Value = -14402.5

$ws.Range("H126").Value = 3098908.8
$ws.Range("J126").Value = 7355287
$ws.Range("L126").Value = 22065861
$ws.Range("N126").Value = -22070801

$ws.Range("H132").Value = 5436.143
$ws.Range("I132").Value = 4138.625
$ws.Range("K132").Value = 12415.875
$ws.Range("M132").Value = -9885.875

$ws.Range("H135").Value = 44903.06
$ws.Range("J135").Value = 44903.06
$ws.Range("L135").Value = 44903.06
$ws.Range("N135").Value = -55043.06

$ws.Range("H141").Value = 65422.668
$ws.Range("J141").Value = 65422.668
$ws.Range("L141").Value = 65422.668
$ws.Range("N141").Value = -75782.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4658.263
$ws.Range("I132").Value = 5007.4287
$ws.Range("J132").Value = 3680.6
$ws.Range("K132").Value = 15022.2861
$ws.Range("L132").Value = 11041.8
$ws.Range("M132").Value = -12492.2861
$ws.Range("N132").Value = -16101.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6529.129
$ws.Range("I132").Value = 4473.3335
$ws.Range("J132").Value = 8456.4375
$ws.Range("K132").Value = 13420.0005
$ws.Range("L132").Value = 25369.3125
$ws.Range("M132").Value = -10890.0005
$ws.Range("N132").Value = -30429.3125

$ws.Range("H136").Value = 20761.482
$ws.Range("I136").Value = 48427.715
$ws.Range("J136").Value = 5059.027
$ws.Range("K136").Value = 145283.145
$ws.Range("L136").Value = 15177.081
$ws.Range("M136").Value = -142733.145
$ws.Range("N136").Value = -20277.081

$ws.Range("H140").Value = 62500
$ws.Range("J140").Value = 62500
$ws.Range("L140").Value = 62500
$ws.Range("N140").Value = -72860

$ws.Range("H141").Value = 58193
$ws.Range("J141").Value = 58193
$ws.Range("L141").Value = 58193
$ws.Range("N141").Value = -68553
